$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-25 04:48:06"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "61%"
$ws.Range("H3").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("E3").Value = "2026-02-25 04:48:08"
$ws.Range("E4").Value = "2026-02-25 04:48:10"
$ws.Range("J4").Value = "1019.5 hPa"
$ws.Range("K4").Value = "-0.1 MJ/m2"
$ws.Range("L4").Value = "5.0 km/h - 102º 4:12 TU"
$ws.Range("O4").Value = "2.8 °C"
$ws.Range("E5").Value = "2026-02-25 04:48:13"
$ws.Range("E6").Value = "2026-02-25 04:48:15"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "95%"
$ws.Range("H3").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "1019.0 hPa"
$ws.Range("N6").Value = "7.1 °C 4:28 TU"
$ws.Range("O6").Value = "8.7 °C"
$ws.Range("E7").Value = "2026-02-25 04:48:17"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "94%"
$ws.Range("H3").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("J7").Value = "1018.6 hPa"
$ws.Range("N7").Value = "10.0 °C 4:23 TU"
$ws.Range("O7").Value = "11.0 °C"
$ws.Range("E8").Value = "2026-02-25 04:48:20"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "47%"
$ws.Range("H3").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("N8").Value = "13.5 °C 4:29 TU"
$ws.Range("O8").Value = "15.2 °C"
$ws.Range("E9").Value = "2026-02-25 04:48:22"
$ws.Range("O9").Value = "5.4 °C"
$ws.Range("E10").Value = "2026-02-25 04:48:24"
$ws.Range("N10").Value = "2.5 °C 4:03 TU"
$ws.Range("O10").Value = "4.0 °C"
$ws.Range("E11").Value = "2026-02-25 04:48:27"
$ws.Range("N11").Value = "2.1 °C 4:00 TU"
$ws.Range("O11").Value = "3.1 °C"
$ws.Range("E12").Value = "2026-02-25 04:48:29"
$ws.Range("M12").Value = "6.8 °C 4:14 TU"
$ws.Range("E13").Value = "2026-02-25 04:48:31"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "93%"
$ws.Range("H3").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("O13").Value = "-1.4 °C"
$ws.Range("E14").Value = "2026-02-25 04:48:33"
$ws.Range("N14").Value = "3.4 °C 4:14 TU"
$ws.Range("O14").Value = "5.6 °C"
$ws.Range("E15").Value = "2026-02-25 04:48:36"
$ws.Range("O15").Value = "5.6 °C"
$ws.Range("E16").Value = "2026-02-25 04:48:38"
$ws.Range("K16").Value = "-0.1 MJ/m2"
$ws.Range("E17").Value = "2026-02-25 04:48:40"
$ws.Range("E18").Value = "2026-02-25 04:48:42"
$ws.Range("J18").Value = "1019.4 hPa"
$ws.Range("N18").Value = "5.0 °C 4:03 TU"
$ws.Range("O18").Value = "6.4 °C"
$ws.Range("E19").Value = "2026-02-25 04:48:45"
$ws.Range("O19").Value = "9.7 °C"
$ws.Range("E20").Value = "2026-02-25 04:48:47"
$ws.Range("M20").Value = "3.4 °C 4:29 TU"
$ws.Range("E21").Value = "2026-02-25 04:48:49"
$ws.Range("J21").Value = "1023.4 hPa"
$ws.Range("O21").Value = "3.8 °C"
$ws.Range("E22").Value = "2026-02-25 04:48:51"
$ws.Range("E23").Value = "2026-02-25 04:48:54"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "27%"
$ws.Range("H3").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("K23").Value = "-0.1 MJ/m2"
$ws.Range("L23").Value = "11.9 km/h - 17º 4:11 TU"
$ws.Range("O23").Value = "3.3 °C"
$ws.Range("E24").Value = "2026-02-25 04:48:56"
$ws.Range("L24").Value = "7.2 km/h - 45º 4:17 TU"
$ws.Range("O24").Value = "3.8 °C"
$ws.Range("E25").Value = "2026-02-25 04:48:58"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "25%"
$ws.Range("H3").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("K25").Value = "-0.1 MJ/m2"
$ws.Range("O25").Value = "3.2 °C"
$ws.Range("E26").Value = "2026-02-25 04:49:00"
$ws.Range("J26").Value = "1018.9 hPa"
$ws.Range("N26").Value = "8.2 °C 4:16 TU"
$ws.Range("O26").Value = "9.3 °C"
$ws.Range("E27").Value = "2026-02-25 04:49:03"
$ws.Range("O27").Value = "4.0 °C"
$ws.Range("E28").Value = "2026-02-25 04:49:05"
$ws.Range("J28").Value = "1020.4 hPa"
$ws.Range("N28").Value = "1.9 °C 4:04 TU"
$ws.Range("O28").Value = "3.9 °C"
$ws.Range("E29").Value = "2026-02-25 04:49:07"
$ws.Range("M29").Value = "10.3 °C 4:29 TU"
$ws.Range("O29").Value = "9.0 °C"
$ws.Range("E30").Value = "2026-02-25 04:49:09"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "99%"
$ws.Range("H3").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("J30").Value = "1019.3 hPa"
$ws.Range("N30").Value = "6.7 °C 4:24 TU"
$ws.Range("O30").Value = "7.7 °C"
$ws.Range("E31").Value = "2026-02-25 04:49:12"
$ws.Range("J31").Value = "1018.5 hPa"
$ws.Range("E32").Value = "2026-02-25 04:49:14"
$ws.Range("O32").Value = "2.4 °C"
$ws.Range("E33").Value = "2026-02-25 04:49:16"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "66%"
$ws.Range("H3").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("N33").Value = "1.1 °C 4:16 TU"
$ws.Range("O33").Value = "2.4 °C"
$ws.Range("E34").Value = "2026-02-25 04:49:18"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "54%"
$ws.Range("H3").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("O34").Value = "1.6 °C"
$ws.Range("E35").Value = "2026-02-25 04:49:20"
$ws.Range("J35").Value = "1019.7 hPa"
$ws.Range("N35").Value = "8.6 °C 4:27 TU"
$ws.Range("O35").Value = "10.0 °C"
$ws.Range("E36").Value = "2026-02-25 04:49:23"
$ws.Range("L36").Value = "37.8 km/h - 12º 4:20 TU"
$ws.Range("M36").Value = "11.9 °C 4:23 TU"
$ws.Range("O36").Value = "8.5 °C"
$ws.Range("E37").Value = "2026-02-25 04:49:25"
$ws.Range("J37").Value = "1024.1 hPa"
$ws.Range("L37").Value = "4.7 km/h - 176º 4:10 TU"
$ws.Range("N37").Value = "0.2 °C 4:06 TU"
$ws.Range("E38").Value = "2026-02-25 04:49:27"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "97%"
$ws.Range("H3").Copy()
$ws.Range("H38").PasteSpecial(-4122)
$ws.Range("E39").Value = "2026-02-25 04:49:30"
$ws.Range("K39").Value = "-0.1 MJ/m2"
$ws.Range("E40").Value = "2026-02-25 04:49:32"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "95%"
$ws.Range("H3").Copy()
$ws.Range("H40").PasteSpecial(-4122)
$ws.Range("J40").Value = "1024.3 hPa"
$ws.Range("N40").Value = "0.5 °C 4:12 TU"
$ws.Range("O40").Value = "1.6 °C"
$ws.Range("E41").Value = "2026-02-25 04:49:34"
$ws.Range("E42").Value = "2026-02-25 04:49:37"
$ws.Range("M42").Value = "10.9 °C 4:24 TU"
$ws.Range("O42").Value = "7.9 °C"
$ws.Range("E43").Value = "2026-02-25 04:49:39"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "92%"
$ws.Range("H3").Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("O43").Value = "3.8 °C"
$ws.Range("E44").Value = "2026-02-25 04:49:41"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "44%"
$ws.Range("H3").Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("O44").Value = "-0.3 °C"
$ws.Range("E45").Value = "2026-02-25 04:49:43"
$ws.Range("O45").Value = "6.6 °C"
$ws.Range("E46").Value = "2026-02-25 04:49:46"
$ws.Range("N46").Value = "2.3 °C 4:08 TU"
$ws.Range("O46").Value = "3.7 °C"
$excel.CutCopyMode = $false
